$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I2").Value = 5.6
$ws.Range("K2").Value = 2.18
$ws.Range("P2").Value = 3.25
$ws.Range("Q2").Value = 1.78
$ws.Range("R2").Value = 1.93
$ws.Range("U2").Value = 1.8
$ws.Range("W2").Value = 6.7
$ws.Range("X2").Value = 7.2
$ws.Range("AH2").Value = 16
$ws.Range("AO2").Value = 7.5
$ws.Range("AQ2").Value = 24
$ws.Range("AT2").Value = 2.6
$ws.Range("AW2").Value = 7
$ws.Range("AX2").Value = 32
